# Wind Load Generator - update the Structural Height input and refresh
# the view state to match where the user left off in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Structural Height" (row 5, column G) changes from 60 to 40.
# All dependent formulas (D14/D17, G14, G17, C31:C35, E31:E35, ...) recalc
# automatically.
$ws.Range("G5").Value = 40

# Leave the selection where the user last clicked before saving.
$ws.Range("G6").Select()
